$d = $word.ActiveDocument

# --- 1 & 4: strike-through the two "Alarm of the event" list items -------
# (Title/Location/Time/Alarm/Notes lists appear twice in the doc: once for
# "select to input details..." and once for "select to change details...".)
$alarmTexts = @(
    "The user may select to input details pertaining to the Alarm of the event.",
    "The user may select to change details pertaining to the Alarm of the event."
)

foreach ($text in $alarmTexts) {
    $r = $d.Content
    $found = $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $paraIndex = $r.Paragraphs.First.Index
        $fullPara = $d.Paragraphs.Item($paraIndex).Range
        $fullPara.Font.StrikeThrough = 1
    }
}

# --- 2 & 3: move the "_GoBack" bookmark ------------------------------------
# It currently sits inside "The system shall display all events to the
# user." (between "all" and " events"); it moves to sit between "The
# system" and " shall reveal all identified changes to the calendar."
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$r2 = $d.Content
$found2 = $r2.Find.Execute("The system shall reveal all identified changes to the calendar.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $splitPoint = $r2.Start + [string]"The system".Length
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
